# Three MAGs (even_MAG-GUT16183.fa, even_MAG-GUT85125.fa, even_MAG-GUT86439.fa)
# were dropped from the results table. Delete their rows (originally rows
# 3, 7 and 9) bottom-up so earlier indices stay valid while deleting,
# letting the remaining rows shift up and the sheet dimension shrink from
# A1:X11 to A1:X8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(3).Delete()
